# Auto-generated Excel COM-interop script
# Applies scheduled market-data updates to the Behemoth_Profits workbook
# (Leve profit sheets: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 60009.5
$ws.Range("I21").Value = 60009.5
$ws.Range("K21").Value = 60009.5
$ws.Range("M21").Value = -59541.5

$ws.Range("H23").Value = 60009.5
$ws.Range("I23").Value = 60009.5
$ws.Range("K23").Value = 60009.5
$ws.Range("M23").Value = -59775.5

$ws.Range("H29").Value = 3823.5715
$ws.Range("I29").Value = 1591.6666
$ws.Range("J29").Value = 5497.5
$ws.Range("K29").Value = 4774.9998
$ws.Range("L29").Value = 16492.5
$ws.Range("M29").Value = -4493.9998
$ws.Range("N29").Value = -17054.5

$ws.Range("H38").Value = 6599.4
$ws.Range("I38").Value = 3749.25
$ws.Range("J38").Value = 18000
$ws.Range("K38").Value = 11247.75
$ws.Range("L38").Value = 54000
$ws.Range("M38").Value = -10875.75
$ws.Range("N38").Value = -54744

$ws.Range("H40").Value = 3563.1765
$ws.Range("J40").Value = 4041.1177
$ws.Range("L40").Value = 4041.1177
$ws.Range("N40").Value = -4391.1177

$ws.Range("H58").Value = 10358.6
$ws.Range("J58").Value = 24996.5
$ws.Range("L58").Value = 74989.5
$ws.Range("N58").Value = -75289.5

$ws.Range("H64").Value = 4823.5293
$ws.Range("I64").Value = 4000
$ws.Range("J64").Value = 4933.3335
$ws.Range("K64").Value = 4000
$ws.Range("L64").Value = 4933.3335
$ws.Range("M64").Value = -3752
$ws.Range("N64").Value = -5429.3335

$ws.Range("H67").Value = 4823.5293
$ws.Range("I67").Value = 4000
$ws.Range("J67").Value = 4933.3335
$ws.Range("K67").Value = 4000
$ws.Range("L67").Value = 4933.3335
$ws.Range("M67").Value = -3142
$ws.Range("N67").Value = -6649.3335

$ws.Range("H69").Value = 26223.428
$ws.Range("I69").Value = 18999.5
$ws.Range("J69").Value = 29113
$ws.Range("K69").Value = 56998.5
$ws.Range("L69").Value = 87339
$ws.Range("M69").Value = -56124.5
$ws.Range("N69").Value = -89087

$ws.Range("H72").Value = 26223.428
$ws.Range("I72").Value = 18999.5
$ws.Range("J72").Value = 29113
$ws.Range("K72").Value = 170995.5
$ws.Range("L72").Value = 262017
$ws.Range("M72").Value = -166627.5
$ws.Range("N72").Value = -270753

$ws.Range("H120").Value = 114000
$ws.Range("J120").Value = 114000
$ws.Range("L120").Value = 114000
$ws.Range("N120").Value = -123676

$ws.Range("H125").Value = 2736
$ws.Range("I125").Value = 2826
$ws.Range("K125").Value = 25434
$ws.Range("M125").Value = -22974

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 57350.75
$ws.Range("I31").Value = 4249.75
$ws.Range("J31").Value = 110451.75
$ws.Range("K31").Value = 4249.75
$ws.Range("L31").Value = 110451.75
$ws.Range("M31").Value = -3955.75
$ws.Range("N31").Value = -111039.75

$ws.Range("H115").Value = 83183.2
$ws.Range("J115").Value = 83183.2
$ws.Range("L115").Value = 83183.2
$ws.Range("N115").Value = -86317.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 66499.2
$ws.Range("I102").Value = 36021.4
$ws.Range("J102").Value = 96977
$ws.Range("K102").Value = 36021.4
$ws.Range("L102").Value = 96977
$ws.Range("M102").Value = -32776.4
$ws.Range("N102").Value = -103467

$ws.Range("H105").Value = 1570.25
$ws.Range("I105").Value = 1437.4286
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1437.4286
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = 309.5714
$ws.Range("N105").Value = -5994

$ws.Range("H119").Value = 70760.75
$ws.Range("J119").Value = 70760.75
$ws.Range("L119").Value = 70760.75
$ws.Range("N119").Value = -80436.75

$ws.Range("H126").Value = 32666.666
$ws.Range("J126").Value = 32666.666
$ws.Range("L126").Value = 32666.666
$ws.Range("N126").Value = -42546.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 110465.75
$ws.Range("J28").Value = 110465.75
$ws.Range("L28").Value = 110465.75
$ws.Range("N28").Value = -110955.75

$ws.Range("H31").Value = 1068629.2
$ws.Range("J31").Value = 1459352.2
$ws.Range("L31").Value = 1459352.2
$ws.Range("N31").Value = -1459942.2

$ws.Range("H34").Value = 1068629.2
$ws.Range("J34").Value = 1459352.2
$ws.Range("L34").Value = 1459352.2
$ws.Range("N34").Value = -1459756.2

$ws.Range("H94").Value = 5147.2
$ws.Range("J94").Value = 5147.2
$ws.Range("L94").Value = 5147.2
$ws.Range("N94").Value = -6049.2

$ws.Range("H95").Value = 14950
$ws.Range("J95").Value = 14950
$ws.Range("L95").Value = 14950
$ws.Range("N95").Value = -20442

$ws.Range("H108").Value = 76964.5
$ws.Range("J108").Value = 76964.5
$ws.Range("L108").Value = 76964.5
$ws.Range("N108").Value = -84644.5

$ws.Range("H122").Value = 3426.4736
$ws.Range("I122").Value = 1616.5834
$ws.Range("J122").Value = 6529.143
$ws.Range("K122").Value = 4849.7502
$ws.Range("L122").Value = 19587.429
$ws.Range("M122").Value = -2399.7502
$ws.Range("N122").Value = -24487.429

$ws.Range("H129").Value = 61212.5
$ws.Range("I129").Value = 50000
$ws.Range("J129").Value = 64950
$ws.Range("K129").Value = 50000
$ws.Range("L129").Value = 64950
$ws.Range("M129").Value = -45000
$ws.Range("N129").Value = -74950

$ws.Range("H134").Value = 718883.9
$ws.Range("I134").Value = 1430055.8
$ws.Range("K134").Value = 4290167.4
$ws.Range("M134").Value = -4287632.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 380805.8
$ws.Range("J12").Value = 559351.5600000001
$ws.Range("L12").Value = 1678054.68
$ws.Range("N12").Value = -1678400.68

$ws.Range("H92").Value = 2501730.8
$ws.Range("I92").Value = 2501730.8
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 7505192.399999999
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -7503944.399999999
$ws.Range("N92").ClearContents()

$ws.Range("H112").Value = 4031.3333
$ws.Range("I112").Value = 2437.8
$ws.Range("J112").Value = 11999
$ws.Range("K112").Value = 7313.400000000001
$ws.Range("L112").Value = 35997
$ws.Range("M112").Value = -6205.400000000001
$ws.Range("N112").Value = -38213

$ws.Range("H129").Value = 25718666
$ws.Range("I129").Value = 1063.8
$ws.Range("J129").Value = 41792170
$ws.Range("K129").Value = 3191.4
$ws.Range("L129").Value = 125376510
$ws.Range("M129").Value = 1808.6
$ws.Range("N129").Value = -125386510

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 72093.336
$ws.Range("J42").Value = 72093.336
$ws.Range("L42").Value = 72093.336
$ws.Range("N42").Value = -73063.336

$ws.Range("H44").Value = 74345.336
$ws.Range("I44").Value = 74504
$ws.Range("K44").Value = 74504
$ws.Range("M44").Value = -73908

$ws.Range("H45").Value = 85326
$ws.Range("J45").Value = 85326
$ws.Range("L45").Value = 85326
$ws.Range("N45").Value = -86444

$ws.Range("H80").Value = 21445.092
$ws.Range("I80").Value = 12270.857
$ws.Range("J80").Value = 37500
$ws.Range("K80").Value = 12270.857
$ws.Range("L80").Value = 37500
$ws.Range("M80").Value = -11272.857
$ws.Range("N80").Value = -39496

$ws.Range("H83").Value = 21445.092
$ws.Range("I83").Value = 12270.857
$ws.Range("J83").Value = 37500
$ws.Range("K83").Value = 61354.285
$ws.Range("L83").Value = 187500
$ws.Range("M83").Value = -56362.285
$ws.Range("N83").Value = -197484

$ws.Range("H104").Value = 79311.664
$ws.Range("J104").Value = 79311.664
$ws.Range("L104").Value = 79311.664
$ws.Range("N104").Value = -86299.664

$ws.Range("H115").Value = 72093.336
$ws.Range("J115").Value = 72093.336
$ws.Range("L115").Value = 72093.336
$ws.Range("N115").Value = -74443.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2443.1428
$ws.Range("I16").Value = 2443.1428
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2443.1428
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2273.1428
$ws.Range("N16").ClearContents()

$ws.Range("H82").Value = 950
$ws.Range("I82").Value = 950
$ws.Range("K82").Value = 950
$ws.Range("M82").Value = -589

$ws.Range("H85").Value = 950
$ws.Range("I85").Value = 950
$ws.Range("K85").Value = 950
$ws.Range("M85").Value = 298

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H141").Value = 69499.5
$ws.Range("J141").Value = 69499.5
$ws.Range("L141").Value = 69499.5
$ws.Range("N141").Value = -79859.5
